$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 20:35"

# Row 4
$ws.Range("B4").Value = 1537769
$ws.Range("C4").Value = 10105
$ws.Range("D4").Value = 351348
$ws.Range("E4").Value = 1094969
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 474
$ws.Range("H4").Value = 91452

# Row 11
$ws.Range("B11").Value = 177213
$ws.Range("C11").Value = 562
$ws.Range("D11").Value = 154600
$ws.Range("E11").Value = 14516
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 48
$ws.Range("H11").Value = 8097

# Row 114
$ws.Range("B114").Value = 796
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 652
$ws.Range("E114").Value = 93
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 51

# Row 119
$ws.Range("A119").Value = "Guinea Ecuatorial"
$ws.Range("B119").Value = 719
$ws.Range("C119").Value = 125
$ws.Range("D119").Value = 22
$ws.Range("E119").Value = 690
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 7

# Row 120
$ws.Range("A120").Value = "Crucero"
$ws.Range("B120").Value = 712
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 651
$ws.Range("E120").Value = 48
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 13

# Row 121
$ws.Range("A121").Value = "Georgia"
$ws.Range("B121").Value = 701
$ws.Range("C121").Value = 6
$ws.Range("D121").Value = 432
$ws.Range("E121").Value = 257
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 12

# Row 122
$ws.Range("A122").Value = "San Marino"
$ws.Range("B122").Value = 654
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 201
$ws.Range("E122").Value = 412
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 41

# Row 123
$ws.Range("A123").Value = "Jordania"
$ws.Range("B123").Value = 613
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 408
$ws.Range("E123").Value = 196
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 9

# Row 133
$ws.Range("B133").Value = 412
$ws.Range("C133").Value = 21
$ws.Range("D133").Value = 110
$ws.Range("E133").Value = 287
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 15

# Row 138
$ws.Range("B138").Value = 335
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 296
$ws.Range("E138").Value = 15
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 24

# Row 146
$ws.Range("B146").Value = 297
$ws.Range("C146").Value = 5
$ws.Range("D146").Value = 203
$ws.Range("E146").Value = 94
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 0

# Row 189
$ws.Range("A189").Value = "Gambia"
$ws.Range("B189").Value = 24
$ws.Range("C189").Value = 1
$ws.Range("D189").Value = 13
$ws.Range("E189").Value = 10
$ws.Range("F189").Value = 0
$ws.Range("G189").Value = 0
$ws.Range("H189").Value = 1

# Row 190
$ws.Range("A190").Value = "Timor Oriental"
$ws.Range("B190").Value = 24
$ws.Range("C190").Value = 0
$ws.Range("D190").Value = 24
$ws.Range("E190").Value = 0
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 0

# Row 196
$ws.Range("A196").Value = "Santa Lucia"
$ws.Range("B196").Value = 18
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 18
$ws.Range("E196").Value = 0
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 0

# Row 197
$ws.Range("A197").Value = "Nueva Caledonia"
$ws.Range("B197").Value = 18
$ws.Range("C197").Value = 0
$ws.Range("D197").Value = 18
$ws.Range("E197").Value = 0
$ws.Range("F197").Value = 0
$ws.Range("G197").Value = 0
$ws.Range("H197").Value = 0

# Row 209
$ws.Range("A209").Value = "Montserrat"
$ws.Range("B209").Value = 11
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 10
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 1

# Row 210
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("B210").Value = 11
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 11
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# Row 211
$ws.Range("A211").Value = "Seychelles"
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 11
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0

# Row 215
$ws.Range("A215").Value = "Sahara Occidental"
$ws.Range("B215").Value = 6
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 6
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0

# Row 216
$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("B216").Value = 6
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 6
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0
